# Edit Attendance Document to include new data
#
# Inserts two new columns ("Attended" / "Length") between the existing
# "TO-DO" and "Links & Docs" columns, fills in sample data for the first
# meeting row, rewords the meeting summary, and re-wires the hyperlink
# that used to live on the "Links & Docs" cell (which shifts two columns
# to the right).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Make room for the two new columns (D:E). This shifts the old
#    "Links & Docs" column (and its hyperlink target cell) from D to F.
# ---------------------------------------------------------------------
$ws.Columns("D:E").Insert()

# The hyperlink that used to sit on D2 keeps its old anchor after the
# column insert, so drop it here - it gets re-created on the correct
# cell (F2) below.
$ws.Range("D2").Hyperlinks.Delete()

# ---------------------------------------------------------------------
# 2. New "Attended" column.
# ---------------------------------------------------------------------
$ws.Range("D1").Value = "Attended"
$ws.Range("D1").Font.Bold = $true

$ws.Range("D2").Value = "07/07"
$ws.Range("D2").Style = "Good"

# ---------------------------------------------------------------------
# 3. New "Length" column (meeting duration, stored as a time value).
# ---------------------------------------------------------------------
$ws.Range("E1").Value = "Length"
$ws.Range("E1").Font.Bold = $true

$ws.Range("E2").Value = 0.0625
$ws.Range("E2").NumberFormat = "[$-F400]h:mm:ss AM/PM"

# ---------------------------------------------------------------------
# 4. Reword the first-meeting summary and let it wrap in its taller row.
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "First meeting with the team. Discussed the list of projects and shortlisted our top 10 choices.  We also procured our more desirable top 3 that we would like to bid for."
$ws.Range("B2").WrapText = $true
$ws.Rows(2).RowHeight = 52

# ---------------------------------------------------------------------
# 5. Re-create the hyperlink on the "Links & Docs" cell, now at F2.
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("F2"), "https://imgur.com/a/1QSYfvk")
$ws.Range("F2").Font.Underline = $true
$ws.Range("F2").Font.Color = 16711680
